# Updated DEU model - 2025-08-31 20:37
$wb = $excel.ActiveWorkbook

# --- Sheet "ev_charging_uc": reorder the comma-separated hour lists in C13/C14 ---
$wsEv = $wb.Worksheets.Item("ev_charging_uc")
$wsEv.Range("C13").Value = "S1aH3,S1aH5,S2aH4,S2aH5,S1aH6,S1aH4,S2aH2,S2aH6,S3aH5,S3aH3,S1aH2,S2aH3,S3aH6,S3aH4,S3aH2"
$wsEv.Range("C14").Value = "S2aH8,S1aH7,S1aH8,S2aH1,S3aH7,S3aH1,S1aH1,S2aH7,S3aH8"

# --- Sheet "re_profiles": recomputed M11:M34 values (now using the 0.000 format already used by column H) ---
$wsRe = $wb.Worksheets.Item("re_profiles")

$mValues = @{
    11 = 0.1938430342782852
    12 = 0.026668071678084892
    13 = 0.02700513099332593
    14 = 0.02738349268768735
    15 = 0.16812362995623994
    16 = 0.05835165936916387
    17 = 0.029097513134079184
    18 = 0.143410013296216
    19 = 0.06072429049925399
    20 = 0.00838341791794081
    21 = 0.008458168104318807
    22 = 0.008449857575589139
    23 = 0.04747623701059971
    24 = 0.016372393639637398
    25 = 0.008213202236041883
    26 = 0.04219049536798028
    27 = 0.036662071087535136
    28 = 0.005330109965263727
    29 = 0.0053335867092306685
    30 = 0.005522686881551576
    31 = 0.031515657407784566
    32 = 0.010408000848768471
    33 = 0.0051164171189268335
    34 = 0.025960862236259676
}

foreach ($r in $mValues.Keys) {
    $cell = $wsRe.Range("M$r")
    $cell.NumberFormat = "0.000"
    $cell.Value = $mValues[$r]
}

# --- Sheet "re_profiles": swap the S1/S2 rows of the Q10:S13 side table ---
$wsRe.Range("Q12").Value = "S2"
$wsRe.Range("R12").Value = 0.16560240645944377
$wsRe.Range("Q13").Value = "S1"
$wsRe.Range("R13").Value = 0.95029816876880036

$wb.Save()
